$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 18897385564.0575
$ws.Range("C2").Value = -28480950180.2576
$ws.Range("D2").Value = 2807219030.24372
$ws.Range("E2").Value = -1165950415.99334

$ws.Range("B3").Value = 100374959448.97363
$ws.Range("C3").Value = 26744879125.375072
$ws.Range("D3").Value = 9701171395.164284
$ws.Range("E3").Value = 2975906159.21861

$ws.Range("B4").Value = 172482580765.8935
$ws.Range("C4").Value = 56026932700.54675
$ws.Range("D4").Value = 17021074916.95565
$ws.Range("E4").Value = 6010088961.88342

$ws.Range("B5").Value = 318919127178.5612
$ws.Range("C5").Value = 110610302258.55699
$ws.Range("D5").Value = 25490533857.065723
$ws.Range("E5").Value = 10171714908.449276

$ws.Range("B6").Value = 1031278835978.76
$ws.Range("C6").Value = 644600752498.008
$ws.Range("D6").Value = 81892558899.4798
$ws.Range("E6").Value = 48047862036.4007
